$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple text fix that reuses an already-existing shared string -----
$ws.Range("B2").Value = "Track"

# --- Remove the old "FreshTyreSet" row (row 8) --------------------------
$ws.Rows(8).Delete()

# --- Make room: 3 new rows for the expanded "Tyre" block, and 2 new rows
#     for the new "Pitstop / Last" entry at the bottom of the sheet. Both
#     inserts are purely structural here; the cell values are filled in
#     below (in the order the original author appears to have typed them).
$ws.Rows("10:12").Insert()
$ws.Rows("70:71").Insert()

# --- Pitstop / Last (typed first by the author) --------------------------
$ws.Range("A70").Value = "Pitstop"
$ws.Range("B70").Value = "Last"
$ws.Range("G70").Value = "[Integer]"
$ws.Range("H70").Value = "#NR# of the last pitstop"

# --- Tyre / Set / Setup / Fresh / Compound block -------------------------
$ws.Range("B10").Value = "Tyre"
$ws.Range("C10").Value = "Set"
$ws.Range("D10").Value = "Setup"
$ws.Range("G10").Value = "[Integer]"
$ws.Range("H10").Value = "Tyre set, when race starts"

$ws.Range("D11").Value = "Fresh"
$ws.Range("G11").Value = "[Integer]"
$ws.Range("H11").Value = "First fresh tyre set"

$ws.Range("C12").Value = "Compound"
$ws.Range("D12").Value = "Setup"
$ws.Range("G12").Value = "[""Wet"" or ""Dry""]"

# row 13 still holds the old "Tyre" label in column B; it no longer belongs there
$ws.Range("B13").Value = ""

# --- Units correction, last edit made by the author ---------------------
$ws.Range("G6").Value = "[Integer] (Lap Time in Milliseconds)"

# --- Update the view: scroll back to top, select G6 ---------------------
$ws.Range("A1").Select()
$ws.Range("G6").Select()
